# "Added numbers to sheets"
#
# Renumber the three sheet tabs, then reproduce the view-state changes:
#   - "Possible Layout" loses the active-tab flag, its selection moves to X25
#   - "Data from MDN" becomes the active tab, keeping its selection at J26

$wb = $excel.ActiveWorkbook

$wsFlex   = $wb.Worksheets.Item(1)   # "Flex Properties To Use"
$wsLayout = $wb.Worksheets.Item(2)   # "Possible Layout"
$wsMdn    = $wb.Worksheets.Item(3)   # "Data from MDN"

$wsFlex.Name   = "(1) Flex Properties To Use"
$wsLayout.Name = "(2) Possible Layout"
$wsMdn.Name    = "(3) Data from MDN"

# "Possible Layout": move the selection, dropping it as the active tab.
$wsLayout.Range("X25").Select()

# "Data from MDN" becomes the active sheet/tab; its own selection (J26) is unchanged.
$wsMdn.Activate()
